# Add the new "OptionDebugOverlay" row (row 56) to the Menu Options sheet,
# matching the CategoryAdvanced / bool option rows already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "CategoryAdvanced"
$ws.Range("B56").Value = 20
$ws.Range("C56").Value = "OptionDebugOverlay"
$ws.Range("D56").Value = "bool"
# Leading apostrophe forces this to stay literal text ("false") instead of
# being auto-converted to the boolean FALSE, matching the Default column's
# text convention used by every other row (e.g. E2 = "true").
$ws.Range("E56").Value = "'false"
$ws.Range("F56").Value = "Enable periodic debug overlay logging with performance stats"
# Value Source column is always blank for every existing row.
$ws.Range("G56").Value = ""
